# Insert a new weekly record at row 172 ("Fruta / hortaliza, semanal").
# Excel's native Insert() shifts rows 172..208 down to 173..209 and copies
# formatting (incl. the date-styled column D) from the row above, matching
# how the original row 172 was formatted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(172).Insert()

# New row 172 is a duplicate of the (now shifted) former row 172 -- which
# landed on row 173 -- with the date/volume/price/kg-price fields updated.
$src = 173
$dst = 172

$ws.Cells.Item($dst, 1).Value2  = $ws.Cells.Item($src, 1).Value2   # A Mercado ID
$ws.Cells.Item($dst, 2).Value2  = $ws.Cells.Item($src, 2).Value2   # B Mercado
$ws.Cells.Item($dst, 3).Value2  = $ws.Cells.Item($src, 3).Value2   # C Región
$ws.Cells.Item($dst, 4).Value2  = 45244                            # D Fecha
$ws.Cells.Item($dst, 5).Value2  = $ws.Cells.Item($src, 5).Value2   # E Codreg
$ws.Cells.Item($dst, 6).Value2  = $ws.Cells.Item($src, 6).Value2   # F Tipo
$ws.Cells.Item($dst, 7).Value2  = $ws.Cells.Item($src, 7).Value2   # G Producto ID
$ws.Cells.Item($dst, 8).Value2  = $ws.Cells.Item($src, 8).Value2   # H Producto
$ws.Cells.Item($dst, 9).Value2  = $ws.Cells.Item($src, 9).Value2   # I Categoría ID
$ws.Cells.Item($dst, 10).Value2 = $ws.Cells.Item($src, 10).Value2  # J Categoría
$ws.Cells.Item($dst, 11).Value2 = $ws.Cells.Item($src, 11).Value2  # K Variedad
$ws.Cells.Item($dst, 12).Value2 = $ws.Cells.Item($src, 12).Value2  # L Calidad
$ws.Cells.Item($dst, 13).Value2 = 50                                # M Volumen
$ws.Cells.Item($dst, 14).Value2 = 13000                             # N Precio mínimo
$ws.Cells.Item($dst, 15).Value2 = 13000                             # O Precio máximo
$ws.Cells.Item($dst, 16).Value2 = 13000                             # P Precio promedio ponderado
$ws.Cells.Item($dst, 17).Value2 = $ws.Cells.Item($src, 17).Value2  # Q Unidad de comercialización
$ws.Cells.Item($dst, 18).Value2 = $ws.Cells.Item($src, 18).Value2  # R Origen
$ws.Cells.Item($dst, 19).Value2 = 3250                              # S Precio $/Kg
$ws.Cells.Item($dst, 20).Value2 = $ws.Cells.Item($src, 20).Value2  # T Kg / unidad
